$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1430179.4
$ws.Range("J17").Value = 1430179.4
$ws.Range("L17").Value = 4290538.199999999
$ws.Range("N17").Value = -4290874.199999999

# Row 33
$ws.Range("H33").Value = 1469.2222
$ws.Range("I33").Value = 210.66667
$ws.Range("J33").Value = 3042.4167
$ws.Range("K33").Value = 210.66667
$ws.Range("L33").Value = 3042.4167
$ws.Range("M33").Value = 18.33332999999999
$ws.Range("N33").Value = -3500.4167

# Row 88
$ws.Range("H88").Value = 898.5625
$ws.Range("I88").Value = 1020.2857
$ws.Range("J88").Value = 864.48
$ws.Range("K88").Value = 1020.2857
$ws.Range("L88").Value = 864.48
$ws.Range("M88").Value = -614.2857
$ws.Range("N88").Value = -1676.48

# Row 91
$ws.Range("H91").Value = 898.5625
$ws.Range("I91").Value = 1020.2857
$ws.Range("J91").Value = 864.48
$ws.Range("K91").Value = 1020.2857
$ws.Range("L91").Value = 864.48
$ws.Range("M91").Value = 383.7143
$ws.Range("N91").Value = -3672.48

# Row 96
$ws.Range("H96").Value = 895.1111
$ws.Range("I96").Value = 732.7222
$ws.Range("J96").Value = 1219.8889
$ws.Range("K96").Value = 2198.1666
$ws.Range("L96").Value = 3659.6667
$ws.Range("M96").Value = -825.1666
$ws.Range("N96").Value = -6405.6667

# Row 106
$ws.Range("H106").Value = 1999.6666
$ws.Range("I106").Value = 1999.6666
$ws.Range("K106").Value = 1999.6666
$ws.Range("M106").Value = -1368.6666

# Row 112
$ws.Range("H112").Value = 6337782
$ws.Range("I112").Value = 779.5
$ws.Range("K112").Value = 2338.5
$ws.Range("M112").Value = -1230.5

# Row 131
$ws.Range("H131").Value = 25897.375
$ws.Range("I131").Value = 25897.375
$ws.Range("K131").Value = 77692.125
$ws.Range("M131").Value = -72652.125

# Row 138
$ws.Range("H138").Value = 3405.69
$ws.Range("I138").Value = 1591.5238
$ws.Range("J138").Value = 3887.9368
$ws.Range("K138").Value = 4774.5714
$ws.Range("L138").Value = 11663.8104
$ws.Range("M138").Value = 365.4286000000002
$ws.Range("N138").Value = -21943.8104

$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 3668.3333
$ws.Range("I21").Value = 3668.3333
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3668.3333
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -3294.3333
$ws.Range("N21").ClearContents()

# Row 32
$ws.Range("H32").Value = 9727.775
$ws.Range("I32").Value = 7565.7085
$ws.Range("K32").Value = 7565.7085
$ws.Range("M32").Value = -7278.7085

# Row 110
$ws.Range("H110").Value = 7379.5
$ws.Range("J110").Value = 11562
$ws.Range("L110").Value = 11562
$ws.Range("N110").Value = -15652

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1344.64
$ws.Range("I94").Value = 1302.8695
$ws.Range("J94").Value = 1825
$ws.Range("K94").Value = 1302.8695
$ws.Range("L94").Value = 1825
$ws.Range("M94").Value = -851.8695
$ws.Range("N94").Value = -2727

# Row 105
$ws.Range("H105").Value = 1487.1177
$ws.Range("I105").Value = 1563.32
$ws.Range("K105").Value = 1563.32
$ws.Range("M105").Value = 183.6800000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 2996.4
$ws.Range("I41").Value = 2996.4
$ws.Range("K41").Value = 2996.4
$ws.Range("M41").Value = -2568.4

# Row 86
$ws.Range("H86").Value = 6015.1113
$ws.Range("I86").Value = 3368.75
$ws.Range("J86").Value = 8132.2
$ws.Range("K86").Value = 3368.75
$ws.Range("L86").Value = 8132.2
$ws.Range("M86").Value = -2245.75
$ws.Range("N86").Value = -10378.2

# Row 89
$ws.Range("H89").Value = 6015.1113
$ws.Range("I89").Value = 3368.75
$ws.Range("J89").Value = 8132.2
$ws.Range("K89").Value = 16843.75
$ws.Range("L89").Value = 40661
$ws.Range("M89").Value = -11227.75
$ws.Range("N89").Value = -51893

# Row 122
$ws.Range("H122").Value = 2744.0625
$ws.Range("J122").Value = 4002.75
$ws.Range("L122").Value = 12008.25
$ws.Range("N122").Value = -16908.25

$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 529
$ws.Range("I117").Value = 529
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1587
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1855
$ws.Range("N117").ClearContents()

# Row 139
$ws.Range("H139").Value = 3921.8462
$ws.Range("I139").Value = 3764.4285
$ws.Range("K139").Value = 11293.2855
$ws.Range("M139").Value = -6153.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4235.8887

# Row 83
$ws.Range("H83").Value = 4235.8887

# Row 102
$ws.Range("H102").Value = 4047.389
$ws.Range("I102").Value = 3994.647
$ws.Range("J102").Value = 4944
$ws.Range("K102").Value = 3994.647
$ws.Range("L102").Value = 4944
$ws.Range("M102").Value = -2372.647
$ws.Range("N102").Value = -8188

# Row 126
$ws.Range("H126").Value = 17110
$ws.Range("I126").Value = 30000
$ws.Range("J126").Value = 4220
$ws.Range("K126").Value = 90000
$ws.Range("L126").Value = 12660
$ws.Range("M126").Value = -87530
$ws.Range("N126").Value = -17600

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5877.2
$ws.Range("I7").Value = 6552.625
$ws.Range("J7").Value = 3175.5
$ws.Range("K7").Value = 6552.625
$ws.Range("L7").Value = 3175.5
$ws.Range("M7").Value = -6440.625
$ws.Range("N7").Value = -3399.5

# Row 40
$ws.Range("H40").Value = 23130
$ws.Range("I40").Value = 31114.428
$ws.Range("J40").Value = 4499.6665
$ws.Range("K40").Value = 31114.428
$ws.Range("L40").Value = 4499.6665
$ws.Range("M40").Value = -30978.428
$ws.Range("N40").Value = -4771.6665

# Row 126
$ws.Range("H126").Value = 5877.2
$ws.Range("I126").Value = 6552.625
$ws.Range("J126").Value = 3175.5
$ws.Range("K126").Value = 19657.875
$ws.Range("L126").Value = 9526.5
$ws.Range("M126").Value = -17187.875
$ws.Range("N126").Value = -14466.5

# Row 132
$ws.Range("H132").Value = 4649.8306
$ws.Range("I132").Value = 3861.4827
$ws.Range("J132").Value = 5411.9
$ws.Range("K132").Value = 11584.4481
$ws.Range("L132").Value = 16235.7
$ws.Range("M132").Value = -9054.4481
$ws.Range("N132").Value = -21295.7

# Row 136
$ws.Range("H136").Value = 3304.7896
$ws.Range("I136").Value = 1956.2333
$ws.Range("K136").Value = 5868.699900000001
$ws.Range("M136").Value = -3318.699900000001

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2207.6
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 2207.6
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 122
$ws.Range("H122").Value = 1847.6875
$ws.Range("I122").Value = 1547.32
$ws.Range("J122").Value = 2920.4285
$ws.Range("K122").Value = 4641.96
$ws.Range("L122").Value = 8761.2855
$ws.Range("M122").Value = -2191.96
$ws.Range("N122").Value = -13661.2855

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 130
$ws.Range("H130").Value = 69000
$ws.Range("I130").Value = 50000
$ws.Range("J130").Value = 88000
$ws.Range("K130").Value = 50000
$ws.Range("L130").Value = 88000
$ws.Range("M130").Value = -44980
$ws.Range("N130").Value = -98040

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132
$ws.Range("H132").Value = 1795.0312
$ws.Range("I132").Value = 1559.0416
